$d = $word.ActiveDocument

# Locate the paragraph containing "LOQ4073: Química Geral II (Requisito fraco)".
# Immediately after it, the document has:
#   - an empty "Normal" paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
# These three paragraphs (including their paragraph marks) must be removed entirely,
# leaving the following empty "Normal" paragraph and the page-break paragraph intact.

$anchor = $d.Content
$anchor.Find.Execute("LOQ4073: Qu" + [char]0x00ED + "mica Geral II (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Find the paragraph index (within the document's Paragraphs collection) of the
# paragraph that contains the anchor text.
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -le $anchor.Start -and $p.Range.End -ge $anchor.End) {
        $targetIndex = $i
        break
    }
    $i = $i + 1
}

if ($targetIndex -eq -1) {
    throw "Could not locate anchor paragraph"
}

$startPara = $d.Paragraphs.Item($targetIndex + 1)
$endPara = $d.Paragraphs.Item($targetIndex + 3)

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
